$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 354, shifting the existing rows 354:373 down to 355:374
$ws.Rows("354").Insert()

# Populate the newly inserted row 354 with the new weekly record
$ws.Cells.Item(354, 1).Value = 3
$ws.Cells.Item(354, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(354, 3).Value = "Coquimbo"
$ws.Cells.Item(354, 4).Value = 44706
$ws.Cells.Item(354, 5).Value = 5
$ws.Cells.Item(354, 6).Value = 100112031
$ws.Cells.Item(354, 7).Value = "Poroto verde"
$ws.Cells.Item(354, 8).Value = "Magnum"
$ws.Cells.Item(354, 9).Value = "Primera"
$ws.Cells.Item(354, 10).Value = 73
$ws.Cells.Item(354, 11).Value = 29000
$ws.Cells.Item(354, 12).Value = 30000
$ws.Cells.Item(354, 13).Value = 29479
$ws.Cells.Item(354, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(354, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(354, 16).Value = 1179
$ws.Cells.Item(354, 17).Value = 25
$ws.Cells.Item(354, 18).Value = "Hortaliza"
